# Refresh cryptocurrency price & 1h-volume data (GitHub Actions scheduled job).
# Column D ("Price") holds plain text (not numbers) in this sheet, so we briefly
# force a text NumberFormat before writing each value -- otherwise Excel would
# auto-coerce numeric-looking strings (e.g. "1.0000") into real numbers and drop
# the significant trailing zeros / thousands-dot formatting used by the source data.
# ClearFormats() afterwards restores the cell to its original (unstyled) state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.874.97"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.728.01"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4811"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2592"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06162"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.726.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06888"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.443"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.87"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.675.30"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9989"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007099"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.949.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.390"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.395"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.045"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.790"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.35%  "

# Rows 28 and 29 swapped: Toncoin now ranks above BitcoinCash.
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.382"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.37"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.927"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07898"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.646"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04574"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.597"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9977"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9204"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.482"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.973"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9992"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.715"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01490"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3817"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.738"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1147"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05357"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.862"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.90"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.236"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.96%  "
